$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill R126:R135 with numeric 0 (previously empty inline-string placeholders)
for ($r = 126; $r -le 135; $r++) {
    $ws.Cells.Item($r, 18).Value = 0
}

# Append new daily rows 136-145 (new trading days scraped for the stock).
# Columns: A Datetime, B Open, C High, D Low, E Close, F Adj Close (left blank),
# G Volume, H Year, I Month, J Day, K Hour, L Minute, M Second, N Week, O isPivot,
# P two_line_structure, Q detect_structure, R backup (left blank).

$ws.Range("A136").Value = 45643
$ws.Range("A136").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B136").Value = 1540
$ws.Range("C136").Value = 1567.699951171875
$ws.Range("D136").Value = 1528.099975585938
$ws.Range("E136").Value = 1563.300048828125
$ws.Range("G136").Value = 2444937
$ws.Range("H136").Value = 2024
$ws.Range("I136").Value = 12
$ws.Range("J136").Value = 17
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 0
$ws.Range("N136").Value = 51
$ws.Range("O136").Value = 0
$ws.Range("P136").Value = 0
$ws.Range("Q136").Value = 0

$ws.Range("A137").Value = 45644
$ws.Range("A137").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B137").Value = 1580
$ws.Range("C137").Value = 1595.849975585938
$ws.Range("D137").Value = 1539.699951171875
$ws.Range("E137").Value = 1579.599975585938
$ws.Range("G137").Value = 2242721
$ws.Range("H137").Value = 2024
$ws.Range("I137").Value = 12
$ws.Range("J137").Value = 18
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 0
$ws.Range("N137").Value = 51
$ws.Range("O137").Value = 0
$ws.Range("P137").Value = 0
$ws.Range("Q137").Value = 0

$ws.Range("A138").Value = 45645
$ws.Range("A138").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B138").Value = 1560
$ws.Range("C138").Value = 1586
$ws.Range("D138").Value = 1554.400024414062
$ws.Range("E138").Value = 1572.349975585938
$ws.Range("G138").Value = 584413
$ws.Range("H138").Value = 2024
$ws.Range("I138").Value = 12
$ws.Range("J138").Value = 19
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 0
$ws.Range("N138").Value = 51
$ws.Range("O138").Value = 0
$ws.Range("P138").Value = 0
$ws.Range("Q138").Value = 0

$ws.Range("A139").Value = 45646
$ws.Range("A139").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B139").Value = 1578.599975585938
$ws.Range("C139").Value = 1588.949951171875
$ws.Range("D139").Value = 1541.300048828125
$ws.Range("E139").Value = 1544.800048828125
$ws.Range("G139").Value = 484886
$ws.Range("H139").Value = 2024
$ws.Range("I139").Value = 12
$ws.Range("J139").Value = 20
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 0
$ws.Range("N139").Value = 51
$ws.Range("O139").Value = 0
$ws.Range("P139").Value = 0
$ws.Range("Q139").Value = 0

$ws.Range("A140").Value = 45649
$ws.Range("A140").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B140").Value = 1530.599975585938
$ws.Range("C140").Value = 1558.349975585938
$ws.Range("D140").Value = 1530.599975585938
$ws.Range("E140").Value = 1545.050048828125
$ws.Range("G140").Value = 375618
$ws.Range("H140").Value = 2024
$ws.Range("I140").Value = 12
$ws.Range("J140").Value = 23
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 0
$ws.Range("N140").Value = 52
$ws.Range("O140").Value = 0
$ws.Range("P140").Value = 0
$ws.Range("Q140").Value = 0

$ws.Range("A141").Value = 45650
$ws.Range("A141").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B141").Value = 1554.599975585938
$ws.Range("C141").Value = 1569.75
$ws.Range("D141").Value = 1541.5
$ws.Range("E141").Value = 1564.75
$ws.Range("G141").Value = 388993
$ws.Range("H141").Value = 2024
$ws.Range("I141").Value = 12
$ws.Range("J141").Value = 24
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 0
$ws.Range("N141").Value = 52
$ws.Range("O141").Value = 0
$ws.Range("P141").Value = 0
$ws.Range("Q141").Value = 0

$ws.Range("A142").Value = 45652
$ws.Range("A142").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B142").Value = 1579.599975585938
$ws.Range("C142").Value = 1579.599975585938
$ws.Range("D142").Value = 1553.050048828125
$ws.Range("E142").Value = 1566.949951171875
$ws.Range("G142").Value = 544734
$ws.Range("H142").Value = 2024
$ws.Range("I142").Value = 12
$ws.Range("J142").Value = 26
$ws.Range("K142").Value = 0
$ws.Range("L142").Value = 0
$ws.Range("M142").Value = 0
$ws.Range("N142").Value = 52
$ws.Range("O142").Value = 0
$ws.Range("P142").Value = 0
$ws.Range("Q142").Value = 0

$ws.Range("A143").Value = 45653
$ws.Range("A143").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B143").Value = 1560.050048828125
$ws.Range("C143").Value = 1589.5
$ws.Range("D143").Value = 1560
$ws.Range("E143").Value = 1579.400024414062
$ws.Range("G143").Value = 356482
$ws.Range("H143").Value = 2024
$ws.Range("I143").Value = 12
$ws.Range("J143").Value = 27
$ws.Range("K143").Value = 0
$ws.Range("L143").Value = 0
$ws.Range("M143").Value = 0
$ws.Range("N143").Value = 52
$ws.Range("O143").Value = 0
$ws.Range("P143").Value = 0
$ws.Range("Q143").Value = 0

$ws.Range("A144").Value = 45656
$ws.Range("A144").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B144").Value = 1584.150024414062
$ws.Range("C144").Value = 1665.349975585938
$ws.Range("D144").Value = 1567.199951171875
$ws.Range("E144").Value = 1619.550048828125
$ws.Range("G144").Value = 3987698
$ws.Range("H144").Value = 2024
$ws.Range("I144").Value = 12
$ws.Range("J144").Value = 30
$ws.Range("K144").Value = 0
$ws.Range("L144").Value = 0
$ws.Range("M144").Value = 0
$ws.Range("N144").Value = 1
$ws.Range("O144").Value = 0
$ws.Range("P144").Value = 0
$ws.Range("Q144").Value = 0

$ws.Range("A145").Value = 45657
$ws.Range("A145").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B145").Value = 1615.050048828125
$ws.Range("C145").Value = 1633.349975585938
$ws.Range("D145").Value = 1601.550048828125
$ws.Range("E145").Value = 1625.449951171875
$ws.Range("G145").Value = 816402
$ws.Range("H145").Value = 2024
$ws.Range("I145").Value = 12
$ws.Range("J145").Value = 31
$ws.Range("K145").Value = 0
$ws.Range("L145").Value = 0
$ws.Range("M145").Value = 0
$ws.Range("N145").Value = 1
$ws.Range("O145").Value = 0
$ws.Range("P145").Value = 0
$ws.Range("Q145").Value = 0

Write-Host "Added rows 136-145; dimension now $($ws.UsedRange.Address())"
